$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.745.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.98%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.525.72"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.69%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.63%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.527"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.524.29"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.81%  "
$ws.Range("E10").Value = "  +0.22%  "
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("E12").Value = "  -3.86%  "
$ws.Range("E13").Value = "  -1.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.983.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000177"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.635.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.519.34"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.70%  "
$ws.Range("E20").Value = "  -5.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "348.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.54%  "
$ws.Range("E22").Value = "  -2.60%  "
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("E24").Value = "  +2.24%  "
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("E27").Value = "  -3.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0983"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "528.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.33"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.07%  "
$ws.Range("E34").Value = "  -2.46%  "
$ws.Range("E35").Value = "  -4.03%  "
$ws.Range("E36").Value = "  -0.15%  "
$ws.Range("E37").Value = "  -2.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.11"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.66"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.03%  "
$ws.Range("E40").Value = "  +0.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.356"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.90%  "
$ws.Range("E42").Value = "  -0.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.11"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.78%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("E45").Value = "  +3.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.46"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.87%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "149.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.95%  "
$ws.Range("E48").Value = "  -3.39%  "
$ws.Range("E49").Value = "  -2.65%  "
$ws.Range("E50").Value = "  +1.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0271"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -10.07%  "
